# Applies the cryptos-list refresh captured in the Sun Nov 3 15:34:09 UTC 2024
# GitHub Actions commit: updates price/volume text cells on Sheet1 (rows 2-51)
# and corrects the NEARProtocol / WrappedeETH row ordering (rows 26-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell, its new text, and whether the text must be
# forced into Text format first (otherwise Excel would silently reinterpret
# a clean-looking numeric string such as "550.80" as the number 550.8).
$updates = @(
    @{Cell='D2'; Value='67.640.59'; ForceText=$false},
    @{Cell='E2'; Value='  -2.44%  '; ForceText=$false},
    @{Cell='D3'; Value='2.424.72'; ForceText=$false},
    @{Cell='E3'; Value='  -2.59%  '; ForceText=$false},
    @{Cell='E4'; Value='  -0.09%  '; ForceText=$false},
    @{Cell='D5'; Value='550.80'; ForceText=$true},
    @{Cell='E5'; Value='  -2.70%  '; ForceText=$false},
    @{Cell='D6'; Value='158.66'; ForceText=$true},
    @{Cell='E7'; Value='  -0.03%  '; ForceText=$false},
    @{Cell='E8'; Value='  -3.02%  '; ForceText=$false},
    @{Cell='D9'; Value='2.421.74'; ForceText=$false},
    @{Cell='E9'; Value='  -2.62%  '; ForceText=$false},
    @{Cell='D10'; Value='0.144'; ForceText=$true},
    @{Cell='E10'; Value='  -8.62%  '; ForceText=$false},
    @{Cell='E11'; Value='  -1.71%  '; ForceText=$false},
    @{Cell='D12'; Value='0.331'; ForceText=$true},
    @{Cell='E12'; Value='  -5.92%  '; ForceText=$false},
    @{Cell='E13'; Value='  -4.17%  '; ForceText=$false},
    @{Cell='D14'; Value='2.869.03'; ForceText=$false},
    @{Cell='E14'; Value='  -2.33%  '; ForceText=$false},
    @{Cell='D15'; Value='67.747.81'; ForceText=$false},
    @{Cell='E15'; Value='  -2.11%  '; ForceText=$false},
    @{Cell='D16'; Value='0.0000164'; ForceText=$true},
    @{Cell='E16'; Value='  -6.17%  '; ForceText=$false},
    @{Cell='D17'; Value='22.84'; ForceText=$true},
    @{Cell='E17'; Value='  -5.99%  '; ForceText=$false},
    @{Cell='D18'; Value='2.438.94'; ForceText=$false},
    @{Cell='E18'; Value='  -2.50%  '; ForceText=$false},
    @{Cell='D19'; Value='10.63'; ForceText=$true},
    @{Cell='E19'; Value='  -4.94%  '; ForceText=$false},
    @{Cell='D20'; Value='336.39'; ForceText=$true},
    @{Cell='E20'; Value='  -2.20%  '; ForceText=$false},
    @{Cell='D21'; Value='6.93'; ForceText=$true},
    @{Cell='E21'; Value='  -5.59%  '; ForceText=$false},
    @{Cell='D22'; Value='3.71'; ForceText=$true},
    @{Cell='E22'; Value='  -3.86%  '; ForceText=$false},
    @{Cell='E23'; Value='  -0.25%  '; ForceText=$false},
    @{Cell='D24'; Value='1.81'; ForceText=$true},
    @{Cell='E24'; Value='  -5.32%  '; ForceText=$false},
    @{Cell='D25'; Value='65.73'; ForceText=$true},
    @{Cell='E25'; Value='  -5.22%  '; ForceText=$false},
    @{Cell='B26'; Value='NEARProtocol'; ForceText=$false},
    @{Cell='C26'; Value='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; ForceText=$false},
    @{Cell='D26'; Value='3.60'; ForceText=$true},
    @{Cell='E26'; Value='  -7.21%  '; ForceText=$false},
    @{Cell='B27'; Value='WrappedeETH'; ForceText=$false},
    @{Cell='C27'; Value='https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; ForceText=$false},
    @{Cell='D27'; Value='2.549.30'; ForceText=$false},
    @{Cell='E27'; Value='  -2.49%  '; ForceText=$false},
    @{Cell='D28'; Value='1.00'; ForceText=$true},
    @{Cell='E28'; Value='  +0.12%  '; ForceText=$false},
    @{Cell='D29'; Value='7.95'; ForceText=$true},
    @{Cell='E29'; Value='  -8.13%  '; ForceText=$false},
    @{Cell='D30'; Value='0.0₃0800'; ForceText=$false},
    @{Cell='E30'; Value='  -8.05%  '; ForceText=$false},
    @{Cell='D31'; Value='6.97'; ForceText=$true},
    @{Cell='E31'; Value='  -8.89%  '; ForceText=$false},
    @{Cell='D32'; Value='0.999'; ForceText=$true},
    @{Cell='E32'; Value='  +0.04%  '; ForceText=$false},
    @{Cell='D33'; Value='418.33'; ForceText=$true},
    @{Cell='E33'; Value='  -5.29%  '; ForceText=$false},
    @{Cell='E34'; Value='  -5.83%  '; ForceText=$false},
    @{Cell='D35'; Value='1.11'; ForceText=$true},
    @{Cell='E35'; Value='  -6.39%  '; ForceText=$false},
    @{Cell='D36'; Value='157.48'; ForceText=$true},
    @{Cell='E36'; Value='  +0.91%  '; ForceText=$false},
    @{Cell='D37'; Value='18.96'; ForceText=$true},
    @{Cell='E37'; Value='  -0.41%  '; ForceText=$false},
    @{Cell='E38'; Value='  -0.23%  '; ForceText=$false},
    @{Cell='E39'; Value='  -5.10%  '; ForceText=$false},
    @{Cell='D40'; Value='17.53'; ForceText=$true},
    @{Cell='E40'; Value='  -3.15%  '; ForceText=$false},
    @{Cell='E41'; Value='  -5.51%  '; ForceText=$false},
    @{Cell='D42'; Value='4.25'; ForceText=$true},
    @{Cell='E42'; Value='  -6.98%  '; ForceText=$false},
    @{Cell='D43'; Value='1.43'; ForceText=$true},
    @{Cell='E43'; Value='  -9.75%  '; ForceText=$false},
    @{Cell='E44'; Value='  -0.83%  '; ForceText=$false},
    @{Cell='D45'; Value='131.86'; ForceText=$true},
    @{Cell='E45'; Value='  -4.52%  '; ForceText=$false},
    @{Cell='D46'; Value='1.98'; ForceText=$true},
    @{Cell='E46'; Value='  -8.12%  '; ForceText=$false},
    @{Cell='E47'; Value='  -4.73%  '; ForceText=$false},
    @{Cell='E48'; Value='  -2.80%  '; ForceText=$false},
    @{Cell='D49'; Value='0.468'; ForceText=$true},
    @{Cell='E49'; Value='  -8.45%  '; ForceText=$false},
    @{Cell='D50'; Value='0.549'; ForceText=$true},
    @{Cell='E50'; Value='  -3.66%  '; ForceText=$false},
    @{Cell='D51'; Value='0.0898'; ForceText=$true},
    @{Cell='E51'; Value='  -2.36%  '; ForceText=$false}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
